$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header style (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New numeric data for columns I and J, rows 2-49
$iValues = @(7,6,6,9,6,6,6,8,8,9,8,8,7,5,6,6,8,8,8,7,9,7,8,6,6,4,5,7,2,7,6,7,5,6,6,3,9,6,9,7,8,6,1,1,6,3,4,4)
$jValues = @(8,7,7,9,8,8,8,8,9,9,8,9,8,6,7,6,8,8,8,7,9,8,8,8,7,8,5,8,6,8,7,8,6,6,7,5,9,8,9,7,8,8,2,5,8,5,5,4)

for ($r = 2; $r -le 49; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
